# Remove the "Upcountry District" column from the invoice template.
# That column (header "Upcountry District" / data "{booking:upcountry_city}")
# lives in column H; deleting it shifts "Upcountry Distance" (old I) and
# "Upcountry Amount" (old J) one column to the left, matching the target.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H").Delete()
